# Append 5 new coded-segment rows (145-149) to Sheet1, mirroring the row
# layout/styling already used by the existing data rows (e.g. row 144),
# and register the handful of brand-new shared strings those rows need.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteSpecial constants used below.
$xlPasteFormats = -4122
$xlPasteValues  = -4163

# Donor cells already present in the sheet that hold (as literal text,
# with the exact style we need) the "document id" / "code" strings that
# rows 145-149 reuse. Re-using them via copy/paste-values keeps the
# existing shared-string entry and avoids Excel re-typing the text as a
# number.
$dDonor = @{ 145 = "D87"; 146 = "D27"; 147 = "D66"; 148 = "D100"; 149 = "D100" }
$eDonor = @{ 145 = "E22"; 146 = "E22"; 147 = "E8";  148 = "E7";   149 = "E8" }
$lDonor = "L120"

# Per-row data for the new records.
$rows = @(
    @{ r = 145; F = "1: 3713"; G = "1: 3733"; H = 0;  I = "Clostridium difficile";  IForceText = $false; J = 21; K = 0.093993; M = "11/8/18 14:04:00" },
    @{ r = 146; F = "1: 1647"; G = "1: 1668"; H = 0;  I = "Nocardia transvalensis"; IForceText = $false; J = 22; K = 0.149976; M = "11/8/18 14:05:00" },
    @{ r = 147; F = "1: 1332"; G = "1: 1336"; H = 0;  I = "2007.";                  IForceText = $true;  J = 5;  K = 0.029303; M = "11/12/18 12:31:00" },
    @{ r = 148; F = "3: 1441"; G = "3: 1447"; H = 0;  I = "October";                IForceText = $false; J = 7;  K = 0.018916; M = "11/12/18 12:32:00" },
    @{ r = 149; F = "3: 1449"; G = "3: 1452"; H = 0;  I = "2013";                   IForceText = $true;  J = 4;  K = 0.010809; M = "11/12/18 12:32:00" }
)

# Scratch cell (well outside the used A1:M range) used to force a value
# that looks numeric ("2013", "2007.") to be stored as literal text: we
# give it an existing text-style (copied from C2, numFmtId 49) so no new
# cell style gets registered, type the value in, then copy *just the
# value* over to the destination (which keeps the destination's own
# style untouched).
$ws.Range("C2").Copy()
$ws.Range("ZZ1").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

foreach ($row in $rows) {
    $r = $row.r

    # Clone the full formatting (styles + row height class) of row 144,
    # the last existing data row, onto the new row.
    $ws.Range("A144:M144").Copy()
    $ws.Range("A$r`:M$r").PasteSpecial($xlPasteFormats)
    $ws.Application.CutCopyMode = $false
    $ws.Rows.Item($r).RowHeight = 16

    # A: bullet marker: identical to every other data row.
    $ws.Range("A$r").Value = "$([char]0x25CF)"

    # B, C: left blank (same as every other row - no comment/group set).

    # D: document id text, reusing an existing shared string/style via copy.
    $ws.Range($dDonor[$r]).Copy()
    $ws.Range("D$r").PasteSpecial($xlPasteValues)
    $ws.Application.CutCopyMode = $false

    # E: code text, reusing an existing shared string/style via copy.
    $ws.Range($eDonor[$r]).Copy()
    $ws.Range("E$r").PasteSpecial($xlPasteValues)
    $ws.Application.CutCopyMode = $false

    # F, G: segment begin/end markers (new shared strings).
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G

    # H: weight score.
    $ws.Range("H$r").Value = $row.H

    # I: segment text (new shared string). A couple of these look like
    # bare numbers ("2013", "2007.") and must stay literal text, so we
    # route them through the scratch cell trick above.
    if ($row.IForceText) {
        $ws.Range("ZZ1").Value = $row.I
        $ws.Range("ZZ1").Copy()
        $ws.Range("I$r").PasteSpecial($xlPasteValues)
        $ws.Application.CutCopyMode = $false
    } else {
        $ws.Range("I$r").Value = $row.I
    }

    # J: area.
    $ws.Range("J$r").Value = $row.J

    # K: coverage percentage.
    $ws.Range("K$r").Value = $row.K

    # L: author, reusing the existing "Sonia" shared string/style.
    $ws.Range($lDonor).Copy()
    $ws.Range("L$r").PasteSpecial($xlPasteValues)
    $ws.Application.CutCopyMode = $false

    # M: creation date/time text (new shared string).
    $ws.Range("M$r").Value = $row.M
}

# Clean up the scratch cell so it doesn't extend the sheet's used range.
# (ClearContents+ClearFormats fully de-registers it from the used range;
# Delete/shift would incorrectly ripple into unrelated rows.)
$ws.Range("ZZ1").ClearContents()
$ws.Range("ZZ1").ClearFormats()
